# Update "Linea 141" horarios workbook: refresh scrape timestamp
# (04:44:38 -> 04:57:25, with one straggler 04:17:03 -> 04:57:25 on
# sheet 3 row 8) and append newly-scraped rows, per commit
# "📊 Horarios actualizados Línea 141 - 727".

$wb = $excel.ActiveWorkbook

$NEW_TS = "04:57:25"

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $NEW_TS"
$ws1.Range("A3").Value = "Total filas: 34"

# Refresh the scrape timestamp + recompute "Minutos" for the rows
# that were re-scraped (all still attached to the same NEW_TS batch).
$ws1.Cells.Item(21,1).Value = $NEW_TS
$ws1.Cells.Item(21,4).Value = 19
$ws1.Cells.Item(22,1).Value = $NEW_TS
$ws1.Cells.Item(22,4).Value = 25
$ws1.Cells.Item(23,1).Value = $NEW_TS
$ws1.Cells.Item(23,4).Value = 37
$ws1.Cells.Item(27,1).Value = $NEW_TS
$ws1.Cells.Item(27,4).Value = 49
$ws1.Cells.Item(28,1).Value = $NEW_TS
$ws1.Cells.Item(28,4).Value = 57
$ws1.Cells.Item(29,1).Value = $NEW_TS
$ws1.Cells.Item(29,4).Value = 67
$ws1.Cells.Item(31,1).Value = $NEW_TS
$ws1.Cells.Item(31,4).Value = 74
$ws1.Cells.Item(33,1).Value = $NEW_TS
$ws1.Cells.Item(33,4).Value = 77
$ws1.Cells.Item(34,1).Value = $NEW_TS
$ws1.Cells.Item(34,4).Value = 84
$ws1.Cells.Item(35,1).Value = $NEW_TS
$ws1.Cells.Item(35,4).Value = 90
$ws1.Cells.Item(36,1).Value = $NEW_TS
$ws1.Cells.Item(36,4).Value = 92
$ws1.Cells.Item(37,1).Value = $NEW_TS
$ws1.Cells.Item(37,4).Value = 94

# New rows appended at the bottom (newly scraped arrivals).
$ws1.Cells.Item(38,1).Value = $NEW_TS
$ws1.Cells.Item(38,2).Value = "06:44"
$ws1.Cells.Item(38,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(38,4).Value = 107
$ws1.Cells.Item(38,5).Value = "LP1912"

$ws1.Cells.Item(39,1).Value = $NEW_TS
$ws1.Cells.Item(39,2).Value = "06:46"
$ws1.Cells.Item(39,3).Value = "215C_EL PATO"
$ws1.Cells.Item(39,4).Value = 109
$ws1.Cells.Item(39,5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $NEW_TS"
$ws2.Range("A3").Value = "Total filas: 14"

$ws2.Cells.Item(15,1).Value = $NEW_TS
$ws2.Cells.Item(15,4).Value = 37
$ws2.Cells.Item(17,1).Value = $NEW_TS
$ws2.Cells.Item(17,4).Value = 74

# New row appended at the bottom.
$ws2.Cells.Item(19,1).Value = $NEW_TS
$ws2.Cells.Item(19,2).Value = "06:46"
$ws2.Cells.Item(19,3).Value = "215C_EL PATO"
$ws2.Cells.Item(19,4).Value = 109
$ws2.Cells.Item(19,5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $NEW_TS"
$ws3.Range("A3").Value = "Total filas: 7"

$ws3.Cells.Item(8,1).Value = $NEW_TS
$ws3.Cells.Item(8,4).Value = 47

# A new arrival (row 10) is inserted ahead of the former last row, which
# shifts down to row 11 unchanged; a further new arrival lands on row 12.
$ws3.Rows.Item(10).Insert()

$ws3.Cells.Item(10,1).Value = $NEW_TS
$ws3.Cells.Item(10,2).Value = "06:09"
$ws3.Cells.Item(10,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(10,4).Value = 72
$ws3.Cells.Item(10,5).Value = "L6173"

$ws3.Cells.Item(12,1).Value = $NEW_TS
$ws3.Cells.Item(12,2).Value = "06:33"
$ws3.Cells.Item(12,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(12,4).Value = 96
$ws3.Cells.Item(12,5).Value = "L6203"
